$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.716.00'
$ws.Range('E2').Value = '  +1.45%  '
$ws.Range('D3').Value = '3.042.61'
$ws.Range('E3').Value = '  +3.27%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = "'380.11"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.32%  '
$ws.Range('D6').Value = "'103.47"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.34%  '
$ws.Range('E7').Value = '  +0.98%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = "'0.597"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.98%  '
$ws.Range('D10').Value = "'37.15"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.79%  '
$ws.Range('E11').Value = '  -0.24%  '
$ws.Range('E12').Value = '  +1.81%  '
$ws.Range('D13').Value = '3.527.24'
$ws.Range('E13').Value = '  +3.41%  '
$ws.Range('D14').Value = "'18.62"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.75%  '
$ws.Range('D15').Value = "'7.77"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.46%  '
$ws.Range('D16').Value = '3.040.66'
$ws.Range('E16').Value = '  +3.17%  '
$ws.Range('D17').Value = "'0.984"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.72%  '
$ws.Range('D18').Value = "'10.53"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -11.53%  '
$ws.Range('D19').Value = '51.768.35'
$ws.Range('E19').Value = '  +1.57%  '
$ws.Range('D20').Value = "'3.06"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.64%  '
$ws.Range('D21').Value = "'12.58"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.44%  '
$ws.Range('D22').Value = '0.0₃0965'
$ws.Range('E22').Value = '  +1.52%  '
$ws.Range('D23').Value = "'70.01"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.86%  '
$ws.Range('D24').Value = "'269.07"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.94%  '
$ws.Range('D25').Value = "'3.17"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.16%  '
$ws.Range('D26').Value = "'8.19"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.14%  '
$ws.Range('D27').Value = "'7.57"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +7.15%  '
$ws.Range('E28').Value = '  +6.60%  '
$ws.Range('D29').Value = "'26.36"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.06%  '
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('D31').Value = "'0.109"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.55%  '
$ws.Range('E32').Value = '  +2.55%  '
$ws.Range('D33').Value = "'34.32"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.63%  '
$ws.Range('B34').Value = 'Toncoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D34').Value = "'2.05"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').Value = "'50.44"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').Value = "'0.0453"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.14%  '
$ws.Range('E37').Value = '  -0.21%  '
$ws.Range('D38').Value = "'3.35"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +8.58%  '
$ws.Range('D39').Value = "'0.291"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +11.51%  '
$ws.Range('D40').Value = "'17.13"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.90%  '
$ws.Range('D41').Value = "'1.88"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.88%  '
$ws.Range('D42').Value = "'2.60"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.20%  '
$ws.Range('E43').Value = '  +0.34%  '
$ws.Range('D44').Value = "'127.63"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.81%  '
$ws.Range('D45').Value = "'3.78"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.82%  '
$ws.Range('D46').Value = "'21.93"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.89%  '
$ws.Range('D47').Value = "'2.13"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.86%  '
$ws.Range('E48').Value = '  +2.90%  '
$ws.Range('D49').Value = '2.039.47'
$ws.Range('E49').Value = '  +1.71%  '
$ws.Range('D50').Value = '3.344.71'
$ws.Range('E50').Value = '  +3.16%  '
$ws.Range('D51').Value = "'0.0322"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.92%  '
